$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly generated data points (YCbCr color-space GLCM run) that were
# previously blank rows in column A.
$ws.Range("A24").Value = 0
$ws.Range("A41").Value = 0
$ws.Range("A43").Value = 0
$ws.Range("A54").Value = 0
$ws.Range("A55").Value = 0
$ws.Range("A56").Value = 0

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("A56").Select()
